# Update ticket/sales counts (column F) for duplicated rows found on the
# "展览" (Exhibition) sheet and on the "全部类型" (All types) sheet.
# Values correspond to the same events listed on both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 4132
$wsExhibit.Range("F6").Value = 493
$wsExhibit.Range("F7").Value = 8595
$wsExhibit.Range("F11").Value = 526
$wsExhibit.Range("F12").Value = 60

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 4132
$wsAll.Range("F8").Value = 493
$wsAll.Range("F10").Value = 8595
$wsAll.Range("F16").Value = 526
$wsAll.Range("F17").Value = 60
